$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old rows 3-11 (rows 1-2 will be overwritten below)
$ws.Range("A1:B11").ClearContents()

# New, smaller result table: H3 / Fe3 / Sum
$ws.Range("A1").Value = "H3"
$ws.Range("B1").Value = 3.021

$ws.Range("A2").Value = "Fe3"
$ws.Range("B2").Value = 167.535

$ws.Range("A4").Value = "Sum"
$ws.Range("B4").Value = 170.556

# Add thin borders around the populated cells
$rng = $ws.Range("A1:B2")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

$rng2 = $ws.Range("A4:B4")
$rng2.Borders.LineStyle = 1
$rng2.Borders.Weight = 2
